$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.861.42'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.42%  '
$ws.Range('D3').Value = '''1.886.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.72%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '''0.7457'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.21%  '
$ws.Range('D6').Value = '''242.18'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('D7').Value = '''1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('D8').Value = '''0.3114'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('D9').Value = '''25.36'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.13%  '
$ws.Range('D10').Value = '''0.07110'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range('D11').Value = '''0.08513'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.88%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').Value = '''0.7590'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.52%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '''1.908.88'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.45%  '
$ws.Range('D14').Value = '''5.351'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.95%  '
$ws.Range('D15').Value = '''93.33'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').Value = '''6.135'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.13%  '
$ws.Range('D17').Value = '''29.902.22'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').Value = '''13.71'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.18%  '
$ws.Range('D19').Value = '''242.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').Value = '''0.000007780'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').Value = '''2.162.18'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').Value = '''7.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.22%  '
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '''0.1589'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.08%  '
$ws.Range('D26').Value = '''9.350'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.52%  '
$ws.Range('D27').Value = '''162.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').Value = '''18.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.48%  '
$ws.Range('E29').Value = '  -0.80%  '
$ws.Range('D30').Value = '''1.512'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.39%  '
$ws.Range('D31').Value = '''1.528'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '''4.467'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = '''4.093'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.25%  '
$ws.Range('D34').Value = '''0.05386'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.96%  '
$ws.Range('D35').Value = '''1.235'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.41%  '
$ws.Range('D36').Value = '''0.7433'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.74%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').Value = '''2.713'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.24%  '
$ws.Range('D39').Value = '''0.01931'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').Value = '''2.773'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.90%  '
$ws.Range('D41').Value = '''0.4448'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.70%  '
$ws.Range('D42').Value = '''1.096.83'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.56%  '
$ws.Range('D43').Value = '''6.068'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').Value = '''72.37'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = '''0.8566'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').Value = '''102.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '''7.650'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').Value = '''1.860'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.36%  '
$ws.Range('D50').Value = '''3.051'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.50%  '
$ws.Range('D51').Value = '''2.045.53'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.18%  '
